# Project2Proposal.pptx edit script
# 1) Update the cached "datetimeFigureOut" field text (18/1/2023 -> 22/1/2023)
#    on every slide layout + the slide master's Date placeholder.
# 2) Append a new slide 4 ("Proj requirements") with a bulleted content list.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date field text wherever it appears (11 layouts + master)
# ---------------------------------------------------------------------------
function Update-DateField($holder) {
    if ($holder -eq $null) { return }
    $hf = $holder.HeadersFooters
    if ($hf -eq $null) { return }
    $dt = $hf.DateAndTime
    if ($dt -eq $null) { return }
    try {
        if ($dt.Text -eq "18/1/2023") {
            $dt.Text = "22/1/2023"
        }
    } catch { }
}

Update-DateField($p.SlideMaster)

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateField($layout)
}

# ---------------------------------------------------------------------------
# 2) Add slide 4: "Proj requirements"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title placeholder ------------------------------------------------------
$title = $s4.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Proj requirements"
$title.LanguageID = "en-SG"

# --- Content placeholder ----------------------------------------------------
$bullets = @(
    "3 components",
    "4 props",
    "1 lifting state",
    "2 states",
    "2 setState",
    "2 routes"
)

$content = $s4.Shapes.Item(2).TextFrame.TextRange
$content.Text = ($bullets -join "`r") + "`r"

$paraCount = $content.Paragraphs().Count

for ($i = 1; $i -le $paraCount; $i++) {
    $para = $content.Paragraphs($i, 1)
    $para.ParagraphFormat.Alignment = 1
    $para.ParagraphFormat.Bullet.Font.Name = "Arial"
    $para.ParagraphFormat.Bullet.Character = 8226
    $para.ParagraphFormat.Bullet.Visible = $true
    $para.LanguageID = "en-US"
    $para.Font.Bold = $false
    $para.Font.Italic = $false
    $para.Font.Color.RGB = 0x2F2924
    $para.Font.Name = "-apple-system"
}

# Split paragraph 5 ("2 setState") into two runs: "2 " + "setState"
$p5 = $content.Paragraphs(5, 1)
$run2 = $p5.Characters(3, 8)
$run2.LanguageID = "en-US"
$run2.Font.Bold = $false
$run2.Font.Italic = $false
$run2.Font.Color.RGB = 0x2F2924
$run2.Font.Name = "-apple-system"
